$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71

# A71 must hold the text "01-04-2021" as a shared string, not be auto-converted
# to a date serial number. Briefly force Text format, enter the value, then
# restore the cell's style to match the other rows (no explicit style index).
$ws.Range("A71").NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-04-2021"
$ws.Range("A71").Style = $ws.Range("A70").Style

$ws.Cells.Item($row, 2).Value = 167616
$ws.Cells.Item($row, 3).Value = 851
$ws.Cells.Item($row, 4).Value = 723
$ws.Cells.Item($row, 5).Value = 72380
$ws.Cells.Item($row, 6).Value = 51978
$ws.Cells.Item($row, 7).Value = 343
$ws.Cells.Item($row, 8).Value = 20059
$ws.Cells.Item($row, 9).Value = 93662
$ws.Cells.Item($row, 10).Value = 23512
$ws.Cells.Item($row, 11).Value = 17706
$ws.Cells.Item($row, 12).Value = 850
$ws.Cells.Item($row, 13).Value = 16856
$ws.Cells.Item($row, 14).Value = 41938
$ws.Cells.Item($row, 15).Value = 38034
$ws.Cells.Item($row, 16).Value = 3904
$ws.Cells.Item($row, 17).Value = 10506
